$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# New values for columns B (case number), C (date serial), E (category text), H (count)
# Row layout: Row, B, C, E, H
$rows = @(
    @{ Row = 2;  B = 260; C = 45051; E = "SAI";          H = 1 },
    @{ Row = 3;  B = 202; C = 45055; E = "SAI";          H = 1 },
    @{ Row = 4;  B = 74;  C = 45048; E = "SAI";          H = 1 },
    @{ Row = 5;  B = 264; C = 45054; E = "SAI";          H = 1 },
    @{ Row = 6;  B = 16;  C = 45055; E = "En el hogar";  H = 1 },
    @{ Row = 7;  B = 148; C = 45054; E = "SAI";          H = 1 },
    @{ Row = 8;  B = 106; C = 45050; E = "SAI";          H = 2 },
    @{ Row = 9;  B = 141; C = 45050; E = "En el hogar";  H = 1 },
    @{ Row = 10; B = 176; C = 45048; E = "En el hogar";  H = 2 },
    @{ Row = 11; B = 201; C = 45056; E = "SAI";          H = 1 },
    @{ Row = 12; B = 307; C = 45049; E = "SAI";          H = 1 },
    @{ Row = 13; B = 186; C = 45055; E = "SAI";          H = 1 },
    @{ Row = 14; B = 192; C = 45055; E = "SAI";          H = 1 },
    @{ Row = 15; B = 236; C = 45056; E = "SAI";          H = 1 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}

# Update the selection on the active sheet view to D1:J1048576 with active cell D1
$ws.Range("D1:J1048576").Select()

# Update the window position for the workbook view
$excel.Left = -120
$excel.Top = -120

$wb.Save()
